# ContactTracingReport.xlsx: rebuild Cases + Contacts table bodies
# with a fresh batch of rows, matching the "Make table scrollable (2)" commit.

$wb = $excel.ActiveWorkbook
$wsCases = $wb.Worksheets.Item("Cases")
$wsContacts = $wb.Worksheets.Item("Contacts")

# --- Clear the old data bodies (keep header row 1 intact) ---
$wsCases.Range("A2:J23").ClearContents()
$wsContacts.Range("A2:F3").ClearContents()

# --- Cases: columns A Id | B Test Date | C Added Date | D Postcode |
#            E Traced? | F Dropped times | G Dropped? | H Traced Date |
#            I Symptom date | J Removed date ---
$r = 2
$wsCases.Cells.Item($r, 1).Value = 1
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3700041551
$wsCases.Cells.Item($r, 4).Value = "CB5"
$wsCases.Cells.Item($r, 5).Value = $true
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false
$wsCases.Cells.Item($r, 8).Value = "04/05/2021 08:54:13"

$r = 3
$wsCases.Cells.Item($r, 1).Value = 2
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.370304294
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $true
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false
$wsCases.Cells.Item($r, 8).Value = "11/05/2021 21:45:11"

$r = 4
$wsCases.Cells.Item($r, 1).Value = 3
$wsCases.Cells.Item($r, 2).Value = 44321
$wsCases.Cells.Item($r, 3).Value = 44320.9039171296
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $true
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false
$wsCases.Cells.Item($r, 8).Value = "04/05/2021 21:46:23"

$r = 5
$wsCases.Cells.Item($r, 1).Value = 4
$wsCases.Cells.Item($r, 2).Value = 44321
$wsCases.Cells.Item($r, 3).Value = 44320.904597662
$wsCases.Cells.Item($r, 4).Value = "CB5"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false
$wsCases.Cells.Item($r, 9).Value = "04/05/2021 00:00:00"

$r = 6
$wsCases.Cells.Item($r, 1).Value = 5
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.987027581
$wsCases.Cells.Item($r, 4).Value = "OX2"
$wsCases.Cells.Item($r, 5).Value = $true
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false
$wsCases.Cells.Item($r, 8).Value = "04/05/2021 23:42:46"

$r = 7
$wsCases.Cells.Item($r, 1).Value = 6
$wsCases.Cells.Item($r, 2).Value = 44321
$wsCases.Cells.Item($r, 3).Value = 44320.9045949074
$wsCases.Cells.Item($r, 4).Value = "CB5"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false
$wsCases.Cells.Item($r, 9).Value = "04/05/2021 00:00:00"

$r = 8
$wsCases.Cells.Item($r, 1).Value = 7
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 9
$wsCases.Cells.Item($r, 1).Value = 8
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 10
$wsCases.Cells.Item($r, 1).Value = 9
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 11
$wsCases.Cells.Item($r, 1).Value = 10
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 12
$wsCases.Cells.Item($r, 1).Value = 11
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 13
$wsCases.Cells.Item($r, 1).Value = 12
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 14
$wsCases.Cells.Item($r, 1).Value = 13
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 15
$wsCases.Cells.Item($r, 1).Value = 14
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 16
$wsCases.Cells.Item($r, 1).Value = 15
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 17
$wsCases.Cells.Item($r, 1).Value = 16
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 18
$wsCases.Cells.Item($r, 1).Value = 17
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 19
$wsCases.Cells.Item($r, 1).Value = 18
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 20
$wsCases.Cells.Item($r, 1).Value = 19
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 21
$wsCases.Cells.Item($r, 1).Value = 20
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 22
$wsCases.Cells.Item($r, 1).Value = 21
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 23
$wsCases.Cells.Item($r, 1).Value = 22
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 24
$wsCases.Cells.Item($r, 1).Value = 23
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 25
$wsCases.Cells.Item($r, 1).Value = 24
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 26
$wsCases.Cells.Item($r, 1).Value = 25
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 27
$wsCases.Cells.Item($r, 1).Value = 26
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 28
$wsCases.Cells.Item($r, 1).Value = 27
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 29
$wsCases.Cells.Item($r, 1).Value = 28
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 30
$wsCases.Cells.Item($r, 1).Value = 29
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 31
$wsCases.Cells.Item($r, 1).Value = 30
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 32
$wsCases.Cells.Item($r, 1).Value = 31
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 33
$wsCases.Cells.Item($r, 1).Value = 32
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 34
$wsCases.Cells.Item($r, 1).Value = 33
$wsCases.Cells.Item($r, 2).Value = 44320
$wsCases.Cells.Item($r, 3).Value = 44320.3703009259
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 35
$wsCases.Cells.Item($r, 1).Value = 34
$wsCases.Cells.Item($r, 2).Value = 44328
$wsCases.Cells.Item($r, 3).Value = 44328.375959213
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

$r = 36
$wsCases.Cells.Item($r, 1).Value = 35
$wsCases.Cells.Item($r, 2).Value = 44328
$wsCases.Cells.Item($r, 3).Value = 44328.3769489583
$wsCases.Cells.Item($r, 4).Value = "OX1"
$wsCases.Cells.Item($r, 5).Value = $false
$wsCases.Cells.Item($r, 6).Value = 0
$wsCases.Cells.Item($r, 7).Value = $false

# --- Expand the Cases table to cover the new rows ---
$loCases = $wsCases.ListObjects.Item(1)
$loCases.Resize($wsCases.Range("A1:J36"))

# --- Contacts: columns A Id | B CaseId | C Added Date | D Traced Date |
#               E Contacted date | F Removed date ---
$r = 2
$wsContacts.Cells.Item($r, 1).Value = 1
$wsContacts.Cells.Item($r, 2).Value = 1
$wsContacts.Cells.Item($r, 3).Value = 44320.3708740625
$wsContacts.Cells.Item($r, 4).Value = "04/05/2021 08:54:13"
$wsContacts.Cells.Item($r, 5).Value = "04/05/2021 08:54:06"

$r = 3
$wsContacts.Cells.Item($r, 1).Value = 2
$wsContacts.Cells.Item($r, 2).Value = 3
$wsContacts.Cells.Item($r, 3).Value = 44320.9065409606
$wsContacts.Cells.Item($r, 4).Value = "04/05/2021 21:46:23"
$wsContacts.Cells.Item($r, 5).Value = "04/05/2021 21:45:31"

$r = 4
$wsContacts.Cells.Item($r, 1).Value = 3
$wsContacts.Cells.Item($r, 2).Value = 5
$wsContacts.Cells.Item($r, 3).Value = 44320.9875444792
$wsContacts.Cells.Item($r, 4).Value = "04/05/2021 23:42:46"
$wsContacts.Cells.Item($r, 5).Value = "04/05/2021 23:42:13"

$r = 5
$wsContacts.Cells.Item($r, 1).Value = 4
$wsContacts.Cells.Item($r, 2).Value = 2
$wsContacts.Cells.Item($r, 3).Value = 44327.90628625
$wsContacts.Cells.Item($r, 4).Value = "11/05/2021 21:45:11"
$wsContacts.Cells.Item($r, 5).Value = "11/05/2021 21:45:11"

# --- Expand the Contacts table to cover the new rows ---
$loContacts = $wsContacts.ListObjects.Item(1)
$loContacts.Resize($wsContacts.Range("A1:F5"))

